# Applies the scraped cryptos-list refresh: per-row Price (D) / Volume(1h) (E)
# updates, plus a name/link/price/volume swap for rows 21-22.
#
# Helper: Excel's Range.Value setter auto-converts numeric-looking strings
# (e.g. "235.64") to real numbers, which would silently reformat values like
# "1.110" -> 1.11 and flip the cell from text to numeric. To preserve the
# exact text, stage the value (apostrophe-prefixed, so Excel keeps it literal)
# in a scratch cell, then copy/paste-special "values only" into the target --
# this carries over the literal text without touching the target cell's style.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("Z1")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.058.84'
$ws.Range("E2").Value = '  -1.49%  '

$ws.Range("D3").Value = '1.850.21'
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("E4").Value = '  +0.05%  '

Set-TextValue $ws.Range("D5") '235.64'
$ws.Range("E5").Value = '  +0.86%  '

$ws.Range("E6").Value = '  +0.03%  '

Set-TextValue $ws.Range("D7") '0.4644'
$ws.Range("E7").Value = '  -1.25%  '

Set-TextValue $ws.Range("D8") '0.2777'
$ws.Range("E8").Value = '  +1.17%  '

Set-TextValue $ws.Range("D9") '0.06408'
$ws.Range("E9").Value = '  +1.14%  '

Set-TextValue $ws.Range("D10") '18.23'
$ws.Range("E10").Value = '  +3.38%  '

Set-TextValue $ws.Range("D11") '97.36'
$ws.Range("E11").Value = '  +15.14%  '

$ws.Range("D12").Value = '1.865.50'
$ws.Range("E12").Value = '  +0.80%  '

Set-TextValue $ws.Range("D13") '0.07512'
$ws.Range("E13").Value = '  +1.28%  '

$ws.Range("E14").Value = '  -1.16%  '

Set-TextValue $ws.Range("D15") '0.6247'
$ws.Range("E15").Value = '  -0.40%  '

Set-TextValue $ws.Range("D16") '294.11'
$ws.Range("E16").Value = '  +21.20%  '

$ws.Range("D17").Value = '29.993.14'
$ws.Range("E17").Value = '  -1.61%  '

Set-TextValue $ws.Range("D18") '1.002'
$ws.Range("E18").Value = '  +0.17%  '

Set-TextValue $ws.Range("D19") '12.63'
$ws.Range("E19").Value = '  -0.39%  '

Set-TextValue $ws.Range("D20") '0.000007362'
$ws.Range("E20").Value = '  +0.21%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.078.59'
$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D22") '1.000'
$ws.Range("E22").Value = '  -0.06%  '

Set-TextValue $ws.Range("D23") '4.982'
$ws.Range("E23").Value = '  +0.49%  '

Set-TextValue $ws.Range("D24") '6.076'
$ws.Range("E24").Value = '  +1.68%  '

Set-TextValue $ws.Range("D25") '165.37'
$ws.Range("E25").Value = '  +1.81%  '

Set-TextValue $ws.Range("D26") '9.067'
$ws.Range("E26").Value = '  -2.10%  '

Set-TextValue $ws.Range("D27") '19.19'
$ws.Range("E27").Value = '  +6.31%  '

$ws.Range("E28").Value = '  +2.69%  '

$ws.Range("E29").Value = '  +6.05%  '

$ws.Range("E30").Value = '  -2.87%  '

Set-TextValue $ws.Range("D31") '3.985'
$ws.Range("E31").Value = '  -1.41%  '

Set-TextValue $ws.Range("D32") '3.821'
$ws.Range("E32").Value = '  -0.89%  '

Set-TextValue $ws.Range("D33") '0.04879'
$ws.Range("E33").Value = '  -0.16%  '

Set-TextValue $ws.Range("D34") '0.7258'
$ws.Range("E34").Value = '  +2.95%  '

Set-TextValue $ws.Range("D35") '1.110'
$ws.Range("E35").Value = '  -2.58%  '

Set-TextValue $ws.Range("D36") '2.730'
$ws.Range("E36").Value = '  +0.81%  '

Set-TextValue $ws.Range("D37") '0.01889'
$ws.Range("E37").Value = '  -0.86%  '

Set-TextValue $ws.Range("D38") '2.652'
$ws.Range("E38").Value = '  -1.05%  '

Set-TextValue $ws.Range("D39") '1.964'
$ws.Range("E39").Value = '  -0.61%  '

Set-TextValue $ws.Range("D40") '0.8545'
$ws.Range("E40").Value = '  -2.03%  '

Set-TextValue $ws.Range("D41") '106.01'
$ws.Range("E41").Value = '  +0.68%  '

Set-TextValue $ws.Range("D42") '1.001'
$ws.Range("E42").Value = '  +0.09%  '

Set-TextValue $ws.Range("D43") '5.676'
$ws.Range("E43").Value = '  +3.29%  '

$ws.Range("E44").Value = '  -1.08%  '

Set-TextValue $ws.Range("D45") '65.48'
$ws.Range("E45").Value = '  +4.73%  '

Set-TextValue $ws.Range("D46") '7.032'
$ws.Range("E46").Value = '  -2.79%  '

Set-TextValue $ws.Range("D47") '8.927'
$ws.Range("E47").Value = '  +4.27%  '

Set-TextValue $ws.Range("D48") '0.1184'
$ws.Range("E48").Value = '  -1.39%  '

Set-TextValue $ws.Range("D49") '33.79'
$ws.Range("E49").Value = '  +1.37%  '

Set-TextValue $ws.Range("D50") '0.05542'
$ws.Range("E50").Value = '  +0.21%  '

Set-TextValue $ws.Range("D51") '0.3701'
$ws.Range("E51").Value = '  +0.28%  '
